$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "test"
$ws.Range("B1").Value = "test"
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "test"

$ws.Range("B2").Select()
